$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Poland Ekstraklasa")

# --- Update existing row 243 (add H/I/J and extend/alter K..AC) ---
$ws.Range("H243").Value = 1
$ws.Range("I243").Value = 2
$ws.Range("J243").Value = "A"
$ws.Range("N243").Value = 2.05
$ws.Range("O243").Value = 3.4
$ws.Range("P243").Value = 3.6
$ws.Range("Q243").Value = -0.5
$ws.Range("U243").Value = 1.95
$ws.Range("V243").Value = 1.9
$ws.Range("W243").Value = -1
$ws.Range("X243").Value = -1
$ws.Range("Y243").Value = 2.6
$ws.Range("Z243").Value = -1
$ws.Range("AA243").Value = 0.825
$ws.Range("AB243").Value = 0.95
$ws.Range("AC243").Value = -1

# --- Copy formatting (style) for new rows 244:250 from row 243 template cells ---
$ws.Range("A243").Copy() | Out-Null
$ws.Range("A244:A250").PasteSpecial(-4122) | Out-Null
$ws.Range("E243").Copy() | Out-Null
$ws.Range("E244:E250").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Append new match rows 244-250 ---
# Row 244
$ws.Range("A244").Value = 242
$ws.Range("B244").Value = 6775592
$ws.Range("C244").Value = "Poland Ekstraklasa"
$ws.Range("D244").Value = "Poland Ekstraklasa"
$ws.Range("E244").Value = 45394.54166666666
$ws.Range("F244").Value = "Warta Poznan"
$ws.Range("G244").Value = "Korona Kielce"
$ws.Range("K244").Value = 2.45
$ws.Range("L244").Value = 3.1
$ws.Range("M244").Value = 2.75
$ws.Range("N244").Value = 2.45
$ws.Range("O244").Value = 3.1
$ws.Range("P244").Value = 2.75
$ws.Range("Q244").Value = 0
$ws.Range("R244").Value = 1.8
$ws.Range("S244").Value = 2.05
$ws.Range("T244").Value = 2
$ws.Range("U244").Value = 1.9
$ws.Range("V244").Value = 1.95
$ws.Range("W244").Value = 0
$ws.Range("X244").Value = 0
$ws.Range("Y244").Value = 0
$ws.Range("Z244").Value = 0
$ws.Range("AA244").Value = 0

# Row 245
$ws.Range("A245").Value = 243
$ws.Range("B245").Value = 6774471
$ws.Range("C245").Value = "Poland Ekstraklasa"
$ws.Range("D245").Value = "Poland Ekstraklasa"
$ws.Range("E245").Value = 45394.64583333334
$ws.Range("F245").Value = "Pogon Szczecin"
$ws.Range("G245").Value = "Ruch Chorzow"
$ws.Range("K245").Value = 1.5
$ws.Range("L245").Value = 4.2
$ws.Range("M245").Value = 5.5
$ws.Range("N245").Value = 1.533
$ws.Range("O245").Value = 4
$ws.Range("P245").Value = 5.25
$ws.Range("Q245").Value = -1
$ws.Range("R245").Value = 1.9
$ws.Range("S245").Value = 1.95
$ws.Range("T245").Value = 3
$ws.Range("U245").Value = 2.025
$ws.Range("V245").Value = 1.825
$ws.Range("W245").Value = 0
$ws.Range("X245").Value = 0
$ws.Range("Y245").Value = 0
$ws.Range("Z245").Value = 0
$ws.Range("AA245").Value = 0

# Row 246
$ws.Range("A246").Value = 244
$ws.Range("B246").Value = 6850053
$ws.Range("C246").Value = "Poland Ekstraklasa"
$ws.Range("D246").Value = "Poland Ekstraklasa"
$ws.Range("E246").Value = 45395.41666666666
$ws.Range("F246").Value = "Stal Mielec"
$ws.Range("G246").Value = "Widzew Lodz"
$ws.Range("K246").Value = 2.8
$ws.Range("L246").Value = 3.3
$ws.Range("M246").Value = 2.3
$ws.Range("N246").Value = 2.8
$ws.Range("O246").Value = 3.3
$ws.Range("P246").Value = 2.3
$ws.Range("Q246").Value = 0.25
$ws.Range("R246").Value = 1.775
$ws.Range("S246").Value = 2.1
$ws.Range("T246").Value = 2.5
$ws.Range("U246").Value = 2.05
$ws.Range("V246").Value = 1.8
$ws.Range("W246").Value = 0
$ws.Range("X246").Value = 0
$ws.Range("Y246").Value = 0
$ws.Range("Z246").Value = 0
$ws.Range("AA246").Value = 0

# Row 247
$ws.Range("A247").Value = 245
$ws.Range("B247").Value = 6774875
$ws.Range("C247").Value = "Poland Ekstraklasa"
$ws.Range("D247").Value = "Poland Ekstraklasa"
$ws.Range("E247").Value = 45395.52083333334
$ws.Range("F247").Value = "Puszcza Niepolomice"
$ws.Range("G247").Value = "Lech Poznan"
$ws.Range("K247").Value = 4.333
$ws.Range("L247").Value = 3.5
$ws.Range("M247").Value = 1.727
$ws.Range("N247").Value = 4.5
$ws.Range("O247").Value = 3.5
$ws.Range("P247").Value = 1.7
$ws.Range("Q247").Value = 0.75
$ws.Range("R247").Value = 1.825
$ws.Range("S247").Value = 2.025
$ws.Range("T247").Value = 2.5
$ws.Range("U247").Value = 2.025
$ws.Range("V247").Value = 1.825
$ws.Range("W247").Value = 0
$ws.Range("X247").Value = 0
$ws.Range("Y247").Value = 0
$ws.Range("Z247").Value = 0
$ws.Range("AA247").Value = 0

# Row 248
$ws.Range("A248").Value = 246
$ws.Range("B248").Value = 6775591
$ws.Range("C248").Value = "Poland Ekstraklasa"
$ws.Range("D248").Value = "Poland Ekstraklasa"
$ws.Range("E248").Value = 45395.625
$ws.Range("F248").Value = "Rakow Czestochowa"
$ws.Range("G248").Value = "Legia Warsaw"
$ws.Range("K248").Value = 2.25
$ws.Range("L248").Value = 3.3
$ws.Range("M248").Value = 3.2
$ws.Range("N248").Value = 2.2
$ws.Range("O248").Value = 3.3
$ws.Range("P248").Value = 3.25
$ws.Range("Q248").Value = -0.25
$ws.Range("R248").Value = 1.875
$ws.Range("S248").Value = 1.975
$ws.Range("T248").Value = 2.5
$ws.Range("U248").Value = 2.025
$ws.Range("V248").Value = 1.825
$ws.Range("W248").Value = 0
$ws.Range("X248").Value = 0
$ws.Range("Y248").Value = 0
$ws.Range("Z248").Value = 0
$ws.Range("AA248").Value = 0

# Row 249
$ws.Range("A249").Value = 247
$ws.Range("B249").Value = 6775593
$ws.Range("C249").Value = "Poland Ekstraklasa"
$ws.Range("D249").Value = "Poland Ekstraklasa"
$ws.Range("E249").Value = 45396.3125
$ws.Range("F249").Value = "LKS Lodz"
$ws.Range("G249").Value = "Radomiak Radom"
$ws.Range("K249").Value = 3.5
$ws.Range("L249").Value = 3.4
$ws.Range("M249").Value = 2.05
$ws.Range("N249").Value = 3.5
$ws.Range("O249").Value = 3.4
$ws.Range("P249").Value = 2.05
$ws.Range("Q249").Value = 0.25
$ws.Range("R249").Value = 2.05
$ws.Range("S249").Value = 1.8
$ws.Range("T249").Value = 2.5
$ws.Range("U249").Value = 2
$ws.Range("V249").Value = 1.85
$ws.Range("W249").Value = 0
$ws.Range("X249").Value = 0
$ws.Range("Y249").Value = 0
$ws.Range("Z249").Value = 0
$ws.Range("AA249").Value = 0

# Row 250
$ws.Range("A250").Value = 248
$ws.Range("B250").Value = 6775589
$ws.Range("C250").Value = "Poland Ekstraklasa"
$ws.Range("D250").Value = "Poland Ekstraklasa"
$ws.Range("E250").Value = 45396.41666666666
$ws.Range("F250").Value = "Jagiellonia Bialystok"
$ws.Range("G250").Value = "Cracovia Krakow"
$ws.Range("K250").Value = 1.833
$ws.Range("L250").Value = 3.5
$ws.Range("M250").Value = 3.8
$ws.Range("N250").Value = 1.75
$ws.Range("O250").Value = 3.6
$ws.Range("P250").Value = 4
$ws.Range("Q250").Value = -0.75
$ws.Range("R250").Value = 2.05
$ws.Range("S250").Value = 1.8
$ws.Range("T250").Value = 2.5
$ws.Range("U250").Value = 1.825
$ws.Range("V250").Value = 2.025
$ws.Range("W250").Value = 0
$ws.Range("X250").Value = 0
$ws.Range("Y250").Value = 0
$ws.Range("Z250").Value = 0
$ws.Range("AA250").Value = 0

